$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the order of the two existing "Sway" / "OneNote ... taking too long to open" rows
$ws.Range("A2").Value = "Microsoft Sway ended with an error is not able to open charts"
$ws.Range("A3").Value = "Microsoft OneNote ended with an error is taking too long to open"

# New topic: "Ishan is not able to run Microsoft Excel" - inserted before the Skype row,
# pushing every row below it down by one
$ws.Rows(8).Insert()
$ws.Range("A8").Value = "Ishan is not able to run Microsoft Excel"
$ws.Range("B8").Value = 1

# New topic: "Ishan is now facing issues with Access" - inserted before the
# "Microsoft Excel ended with an error is not able to open charts" row
$ws.Rows(12).Insert()
$ws.Range("A12").Value = "Ishan is now facing issues with Access"
$ws.Range("B12").Value = 1

# New topic: "Hi I am not able to open Excel" - inserted right after the Publisher row
$ws.Rows(17).Insert()
$ws.Range("A17").Value = "Hi I am not able to open Excel"
$ws.Range("B17").Value = 1
